$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 11.5382
$ws.Range("D3").Value = -7.66919999999999
$ws.Range("A4").Value = -21.04569999999999
$ws.Range("C4").Value = -11.03519999999999
$ws.Range("D4").Value = -6.569999999999998
$ws.Range("C5").Value = -14.76900000000001
$ws.Range("A6").Value = -21.42780000000001
$ws.Range("A7").Value = -21.66630000000001
$ws.Range("C8").Value = -12.3364
$ws.Range("D9").Value = -7.322399999999996
$ws.Range("D11").Value = -8.429099999999996
$ws.Range("D14").Value = -6.267299999999995
$ws.Range("A16").Value = -20.28649999999999
$ws.Range("C16").Value = -12.07260000000001
$ws.Range("D18").Value = -8.412599999999991
$ws.Range("A20").Value = -22.74480000000001
$ws.Range("E20").Value = 12.95449999999999
$ws.Range("C22").Value = -11.12899999999999
$ws.Range("D25").Value = -8.506699999999993
